$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Location Id" column (old column H) ---
$ws.Columns("H").Delete()

# --- Move "Regional Client Id" (old column G, now column G after the delete
#     above shifted nothing left of H) to be right after "Sr No." (column B) ---
$ws.Columns("G").Cut()
$ws.Columns("B").Insert()

# --- The "Email" hyperlink cell shifted from I2 to H2 along with the column
#     operations above, but the engine does not auto-update the worksheet's
#     Hyperlinks collection, so the relationship is re-pointed explicitly. ---
$ws.Range("I2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:testcnr@yopmail.com") | Out-Null

# --- Re-select the whole used range (mirrors the "select all / autofit"
#     gesture that produced the refreshed column widths below). ---
$ws.Cells.Select() | Out-Null

# --- Re-apply best-fit-style column widths for the new A:P layout. ---
$ws.Columns("A").ColumnWidth = 6.11
$ws.Columns("B").ColumnWidth = 15.66
$ws.Columns("C").ColumnWidth = 10.11
$ws.Columns("D").ColumnWidth = 10.78
$ws.Columns("E").ColumnWidth = 16.55
$ws.Columns("F").ColumnWidth = 13.11
$ws.Columns("G").ColumnWidth = 11
$ws.Columns("H").ColumnWidth = 19.22
$ws.Columns("I").ColumnWidth = 12.55
$ws.Columns("J").ColumnWidth = 12.55
$ws.Columns("K").ColumnWidth = 12.55
$ws.Columns("L").ColumnWidth = 12.55
$ws.Columns("M").ColumnWidth = 7.89
$ws.Columns("N").ColumnWidth = 14.66
$ws.Columns("O").ColumnWidth = 10.89
$ws.Columns("P").ColumnWidth = 6.44
